$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Book Inventory")

# New book record pulled in by ISBN lookup: "The Hunger Games" by Suzanne Collins.
$title  = "The Hunger Games"
$author = "Suzanne Collins"
# Leading apostrophe keeps the ISBN as text so the leading zero survives.
$isbn   = "'0439023483"

$ws.Cells.Item(2, 1).Value = $title
$ws.Cells.Item(2, 2).Value = $author
$ws.Cells.Item(2, 3).Value = $isbn

$ws.Cells.Item(3, 1).Value = $title
$ws.Cells.Item(3, 2).Value = $author
$ws.Cells.Item(3, 3).Value = $isbn

$ws.Cells.Item(4, 1).Value = $title
$ws.Cells.Item(4, 2).Value = $author
$ws.Cells.Item(4, 3).Value = $isbn

$ws.Range("C1").Select() | Out-Null
